# Implement Funding Source
#
# 1. Insert a new "fundingSource" column immediately to the left of the
#    existing "title" column on every worksheet that models a catalog/product
#    style entity (Service, Product, Antibody, Hybridoma, Protein, NucleicAcid,
#    DetectionKit, Bundle, Virus, Bacterium, Fungus, Protozoan, Viroid, Prion).
# 2. Append a brand new "FundingSource" worksheet at the end of the workbook
#    describing the funding source entity itself.

$wb = $excel.ActiveWorkbook

$targetSheets = @(
    "Service",
    "Product",
    "Antibody",
    "Hybridoma",
    "Protein",
    "NucleicAcid",
    "DetectionKit",
    "Bundle",
    "Virus",
    "Bacterium",
    "Fungus",
    "Protozoan",
    "Viroid",
    "Prion"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count

    $titleCol = -1
    for ($c = 1; $c -le $lastCol; $c++) {
        $val = $ws.Cells.Item(1, $c).Value2
        if ($val -eq "title") {
            $titleCol = $c
            break
        }
    }

    if ($titleCol -gt 0) {
        # Push "title" (and everything after it) one column to the right,
        # then write the new header into the freshly inserted, empty column.
        $ws.Cells.Item(1, $titleCol).EntireColumn.Insert()
        $ws.Cells.Item(1, $titleCol).Value2 = "fundingSource"
    }
}

# Append the new FundingSource sheet as the last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$fundingSourceSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$fundingSourceSheet.Name = "FundingSource"

$fundingSourceHeaders = @(
    "title",
    "description",
    "fundingProgram",
    "grantNumber",
    "funder",
    "fundingPeriodStart",
    "fundingPeriodEnd",
    "eligibilityCriteria",
    "keyword",
    "dateIssued",
    "dateModified",
    "identifier",
    "iri"
)

for ($i = 0; $i -lt $fundingSourceHeaders.Count; $i++) {
    $fundingSourceSheet.Cells.Item(1, $i + 1).Value2 = $fundingSourceHeaders[$i]
}
